$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 6664.6665
$ws.Range("I43").Value = 5992
$ws.Range("J43").Value = 6799.2
$ws.Range("K43").Value = 5992
$ws.Range("L43").Value = 6799.2
$ws.Range("M43").Value = -5923
$ws.Range("N43").Value = -6937.2
# Row 69
$ws.Range("H69").Value = 8278
$ws.Range("I69").Value = 6337.6665
$ws.Range("J69").Value = 9733.25
$ws.Range("K69").Value = 19012.9995
$ws.Range("L69").Value = 29199.75
$ws.Range("M69").Value = -18138.9995
$ws.Range("N69").Value = -30947.75
# Row 72
$ws.Range("H72").Value = 8278
$ws.Range("I72").Value = 6337.6665
$ws.Range("J72").Value = 9733.25
$ws.Range("K72").Value = 57038.9985
$ws.Range("L72").Value = 87599.25
$ws.Range("M72").Value = -52670.9985
$ws.Range("N72").Value = -96335.25
# Row 86
$ws.Range("H86").Value = 20000
$ws.Range("I86").Value = 20000
$ws.Range("K86").Value = 20000
$ws.Range("M86").Value = -18877
# Row 89
$ws.Range("H89").Value = 20000
$ws.Range("I89").Value = 20000
$ws.Range("K89").Value = 100000
$ws.Range("M89").Value = -94384
# Row 93
$ws.Range("H93").Value = 43000
$ws.Range("J93").Value = 43000
$ws.Range("L93").Value = 43000
$ws.Range("N93").Value = -47992
# Row 103
$ws.Range("H103").Value = 598.8
$ws.Range("I103").Value = 502.33334
$ws.Range("K103").Value = 1507.00002
$ws.Range("M103").Value = -921.0000199999999
# Row 106
$ws.Range("H106").Value = 9566
$ws.Range("I106").Value = 5771.909
$ws.Range("J106").Value = 19999.75
$ws.Range("K106").Value = 5771.909
$ws.Range("L106").Value = 19999.75
$ws.Range("M106").Value = -5140.909
$ws.Range("N106").Value = -21261.75
# Row 107
$ws.Range("H107").Value = 5005.722
$ws.Range("I107").Value = 3888.4119
$ws.Range("J107").Value = 24000
$ws.Range("K107").Value = 3888.4119
$ws.Range("L107").Value = 24000
$ws.Range("M107").Value = -1968.4119
$ws.Range("N107").Value = -27840
# Row 138
$ws.Range("H138").Value = 3393.3044
$ws.Range("I138").Value = 3034.125
$ws.Range("J138").Value = 4214.2856
$ws.Range("K138").Value = 9102.375
$ws.Range("L138").Value = 12642.8568
$ws.Range("M138").Value = -3962.375
$ws.Range("N138").Value = -22922.8568

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1145.7894
$ws.Range("I2").Value = 975.4
$ws.Range("J2").Value = 1784.75
$ws.Range("K2").Value = 975.4
$ws.Range("L2").Value = 1784.75
$ws.Range("M2").Value = -862.4
$ws.Range("N2").Value = -2010.75
# Row 19
$ws.Range("H19").Value = 1008
$ws.Range("I19").Value = 1008
$ws.Range("K19").Value = 1008
$ws.Range("M19").Value = -779
# Row 32
$ws.Range("H32").Value = 177270.23
$ws.Range("I32").Value = 203407.36
$ws.Range("J32").Value = 13913.25
$ws.Range("K32").Value = 203407.36
$ws.Range("L32").Value = 13913.25
$ws.Range("M32").Value = -203120.36
$ws.Range("N32").Value = -14487.25
# Row 61
$ws.Range("H61").Value = 2529.8708
$ws.Range("I61").Value = 2172.3845
$ws.Range("K61").Value = 2172.3845
$ws.Range("M61").Value = -1960.3845
# Row 63
$ws.Range("H63").Value = 3239.8
$ws.Range("J63").Value = 3733
$ws.Range("L63").Value = 3733
$ws.Range("N63").Value = -5105
# Row 66
$ws.Range("H66").Value = 3239.8
$ws.Range("J66").Value = 3733
$ws.Range("L66").Value = 18665
$ws.Range("N66").Value = -25529
# Row 116
$ws.Range("H116").Value = 1145.7894
$ws.Range("I116").Value = 975.4
$ws.Range("J116").Value = 1784.75
$ws.Range("K116").Value = 975.4
$ws.Range("L116").Value = 1784.75
$ws.Range("M116").Value = 1318.6
$ws.Range("N116").Value = -6372.75
# Row 122
$ws.Range("H122").Value = 2454.7188
$ws.Range("I122").Value = 1816
$ws.Range("K122").Value = 5448
$ws.Range("M122").Value = -2998
# Row 132
$ws.Range("H132").Value = 5048.88
$ws.Range("I132").Value = 3685.6453
$ws.Range("K132").Value = 11056.9359
$ws.Range("M132").Value = -8526.9359
# Row 134
$ws.Range("H134").Value = 78749.5
$ws.Range("J134").Value = 78749.5
$ws.Range("L134").Value = 78749.5
$ws.Range("N134").Value = -88889.5
# Row 136
$ws.Range("H136").Value = 2529.8708
$ws.Range("I136").Value = 2172.3845
$ws.Range("K136").Value = 6517.1535
$ws.Range("M136").Value = -3967.1535

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1145.7894
$ws.Range("I3").Value = 975.4
$ws.Range("J3").Value = 1784.75
$ws.Range("K3").Value = 975.4
$ws.Range("L3").Value = 1784.75
$ws.Range("M3").Value = -861.4
$ws.Range("N3").Value = -2012.75
# Row 22
$ws.Range("H22").Value = 1708.5
$ws.Range("I22").Value = 234.83333
$ws.Range("K22").Value = 234.83333
$ws.Range("M22").Value = -61.83332999999999
# Row 134
$ws.Range("H134").Value = 4545.591
$ws.Range("I134").Value = 4500.15
$ws.Range("K134").Value = 13500.45
$ws.Range("M134").Value = -10965.45

$ws = $wb.Worksheets.Item("CRP")
# Row 25
$ws.Range("H25").Value = 9254.5
$ws.Range("I25").Value = 9254.5
$ws.Range("K25").Value = 9254.5
$ws.Range("M25").Value = -9080.5
# Row 132
$ws.Range("H132").Value = 3003.8635
$ws.Range("I132").Value = 2774.25
$ws.Range("K132").Value = 8322.75
$ws.Range("M132").Value = -5792.75
# Row 134
$ws.Range("H134").Value = 3349.7368
$ws.Range("I134").Value = 2801.5
$ws.Range("K134").Value = 8404.5
$ws.Range("M134").Value = -5869.5

$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 560
$ws.Range("I98").Value = 654.5
$ws.Range("K98").Value = 1963.5
$ws.Range("M98").Value = -465.5
# Row 99
$ws.Range("H99").Value = 12337.5
$ws.Range("J99").Value = 15000
$ws.Range("L99").Value = 45000
$ws.Range("N99").Value = -49492
# Row 102
$ws.Range("H102").Value = 4000
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
# Row 103
$ws.Range("H103").Value = 2024.6
$ws.Range("J103").Value = 3199.3333
$ws.Range("L103").Value = 9597.999899999999
$ws.Range("N103").Value = -11355.9999
# Row 114
$ws.Range("H114").Value = 842.2727
$ws.Range("I114").Value = 196.14285
$ws.Range("J114").Value = 1973
$ws.Range("K114").Value = 588.4285500000001
$ws.Range("L114").Value = 5919
$ws.Range("M114").Value = 2665.57145
$ws.Range("N114").Value = -12427
# Row 128
$ws.Range("H128").Value = 293115
$ws.Range("I128").Value = 293115
$ws.Range("K128").Value = 879345
$ws.Range("M128").Value = -874365
# Row 131
$ws.Range("H131").Value = 4690.5527
$ws.Range("J131").Value = 4545.0938
$ws.Range("L131").Value = 13635.2814
$ws.Range("N131").Value = -23715.2814
# Row 140
$ws.Range("H140").Value = 9997.5
$ws.Range("I140").Value = 9997.5
$ws.Range("K140").Value = 29992.5
$ws.Range("M140").Value = -24812.5

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 4253
$ws.Range("J21").Value = 4007
$ws.Range("L21").Value = 4007
$ws.Range("N21").Value = -4353
# Row 30
$ws.Range("H30").Value = 4253
$ws.Range("J30").Value = 4007
$ws.Range("L30").Value = 4007
$ws.Range("N30").Value = -4217
# Row 95
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
# Row 102
$ws.Range("H102").Value = 4089.4443
$ws.Range("I102").Value = 4463.125
$ws.Range("K102").Value = 4463.125
$ws.Range("M102").Value = -2841.125
# Row 126
$ws.Range("H126").Value = 2570.3635
$ws.Range("I126").Value = 2484.25
$ws.Range("K126").Value = 7452.75
$ws.Range("M126").Value = -4982.75
# Row 132
$ws.Range("H132").Value = 8878.906000000001
$ws.Range("I132").Value = 6223.407
$ws.Range("K132").Value = 18670.221
$ws.Range("M132").Value = -16140.221

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 3053.16
$ws.Range("I132").Value = 2379.7058
$ws.Range("K132").Value = 7139.117400000001
$ws.Range("M132").Value = -4609.117400000001

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3322
$ws.Range("I132").Value = 2872.5
$ws.Range("K132").Value = 8617.5
$ws.Range("M132").Value = -6087.5
# Row 136
$ws.Range("H136").Value = 3721.6
$ws.Range("I136").Value = 3579.5557
$ws.Range("K136").Value = 10738.6671
$ws.Range("M136").Value = -8188.667099999999
